# Insert two new quarterly columns (D:E) into the FR worksheet, shifting the
# existing quarters from D:K to F:M, then populate the two new columns with
# the latest reported quarter figures and patch a handful of restated totals
# in the shifted columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank columns before column D; existing D:K data/formatting
# shifts right to F:M automatically.
$ws.Columns("D:E").Insert()

# New D:E columns come in unformatted; clone the number formats/styles from
# column F (the old column D, now shifted) so dates/numbers render the same.
# Done per statement block (Income Statement / Balance Sheet / Cash Flow)
# so the blank separator rows between them are left untouched.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = most recent quarter, E = prior
# quarter) with the newly reported figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 105100
$ws.Range("E8").Value2 = 100300
$ws.Range("D9").Value2 = 30400
$ws.Range("E9").Value2 = 28500
$ws.Range("D10").Value2 = 74700
$ws.Range("E10").Value2 = 71800
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 30900
$ws.Range("E15").Value2 = 28600
$ws.Range("D17").Value2 = 67600
$ws.Range("E17").Value2 = 63700
$ws.Range("D18").Value2 = 37500
$ws.Range("E18").Value2 = 36600
$ws.Range("D20").Value2 = 28200
$ws.Range("E20").Value2 = 7900
$ws.Range("D21").Value2 = 97400
$ws.Range("E21").Value2 = 74600
$ws.Range("D22").Value2 = 13800
$ws.Range("E22").Value2 = 13300
$ws.Range("D23").Value2 = 51900
$ws.Range("E23").Value2 = 31200
$ws.Range("D24").Value2 = 0
$ws.Range("E24").Value2 = -300
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 51900
$ws.Range("E26").Value2 = 31500
$ws.Range("D27").Value2 = 50700
$ws.Range("E27").Value2 = 30800
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -28200
$ws.Range("E32").Value2 = -7900
$ws.Range("D33").Value2 = 50700
$ws.Range("E33").Value2 = 30800
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 50700
$ws.Range("E35").Value2 = 30800
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 43100
$ws.Range("E41").Value2 = 52700
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 76300
$ws.Range("E43").Value2 = 74200
$ws.Range("D44").Value2 = 59500
$ws.Range("E44").Value2 = 142300
$ws.Range("D45").Value2 = 101200
$ws.Range("E45").Value2 = 106200
$ws.Range("D46").Value2 = 0
$ws.Range("E46").Value2 = 0
$ws.Range("D47").Value2 = 23300
$ws.Range("E47").Value2 = 23400
$ws.Range("D48").Value2 = 2802400
$ws.Range("E48").Value2 = 2679700
$ws.Range("D49").Value2 = 29700
$ws.Range("E49").Value2 = 27900
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 7300
$ws.Range("E52").Value2 = 19000
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 3142700
$ws.Range("E54").Value2 = 3125300
$ws.Range("D57").Value2 = 78700
$ws.Range("E57").Value2 = 82400
$ws.Range("D58").Value2 = 0
$ws.Range("E58").Value2 = 0
$ws.Range("D59").Value2 = 28800
$ws.Range("E59").Value2 = 28700
$ws.Range("D60").Value2 = 0
$ws.Range("E60").Value2 = 0
$ws.Range("D61").Value2 = 1297800
$ws.Range("E61").Value2 = 1299100
$ws.Range("D62").Value2 = 0
$ws.Range("E62").Value2 = 0
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 1497200
$ws.Range("E66").Value2 = 1496400
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = -490800
$ws.Range("E72").Value2 = -514100
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 1645500
$ws.Range("E76").Value2 = 1628900
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 50700
$ws.Range("E81").Value2 = 30800
$ws.Range("D83").Value2 = 31700
$ws.Range("E83").Value2 = 30100
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 59500
$ws.Range("E89").Value2 = 57200
$ws.Range("D91").Value2 = -66000
$ws.Range("E91").Value2 = -4100
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -50900
$ws.Range("E94").Value2 = -33000
$ws.Range("D96").Value2 = -28000
$ws.Range("E96").Value2 = -28000
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -29900
$ws.Range("E100").Value2 = -30000
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = 0
$ws.Range("D102").Value2 = -21300
$ws.Range("E102").Value2 = -5800

# A handful of shifted cells were restated along with the new data and don't
# match a straight carry-over from the old columns; patch those explicitly.
$ws.Range("I89").Value2 = 54300
$ws.Range("F91").Value2 = -37900
$ws.Range("G91").Value2 = -49900
$ws.Range("H91").Value2 = -41300
$ws.Range("I91").Value2 = -32500
$ws.Range("J91").Value2 = -72100
$ws.Range("H94").Value2 = 82100
$ws.Range("H102").Value2 = 25400
